$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 116; $row++) {
    $ws.Cells.Item($row, 3).Value2 = 45184
}
